$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing 12 data rows (FAPs/MuSCs sending-cluster blocks, rows 2-13) down to
# rows 8-19 to make room for the new "ECs" sending-cluster block at rows 2-7. Copying values
# manually (rather than Rows.Insert()) avoids pulling in a spurious extra cell style.
for ($r = 13; $r -ge 2; $r--) {
    $destRow = $r + 6
    for ($c = 1; $c -le 20; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $destCell = $ws.Cells.Item($destRow, $c)
        $destCell.Value = $srcCell.Value()
    }
}

# Write out every data cell (rows 2-19) with the refreshed TPM-derived values
# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Bmp7"
$ws.Cells.Item(2, 3).Value = "Acvr2b"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.1301303333333333
$ws.Cells.Item(2, 8).Value = 0.390391
$ws.Cells.Item(2, 9).Value = 0.2986126887311651
$ws.Cells.Item(2, 10).Value = 0.2986126887311651
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.978762
$ws.Cells.Item(2, 14).Value = 2.936286
$ws.Cells.Item(2, 15).Value = 0.3819465121442868
$ws.Cells.Item(2, 16).Value = 0.3819465121442868
$ws.Cells.Item(2, 17).Value = 0.127366625314
$ws.Cells.Item(2, 18).Value = 1.146299627826
$ws.Cells.Item(2, 19).Value = 0.1140540749428961
$ws.Cells.Item(2, 20).Value = 0.1140540749428961

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Bmp7"
$ws.Cells.Item(3, 3).Value = "Acvr2b"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.1301303333333333
$ws.Cells.Item(3, 8).Value = 0.390391
$ws.Cells.Item(3, 9).Value = 0.2986126887311651
$ws.Cells.Item(3, 10).Value = 0.2986126887311651
$ws.Cells.Item(3, 11).Value = 2
$ws.Cells.Item(3, 12).Value = 0.6666666666666666
$ws.Cells.Item(3, 13).Value = 0.1041576666666667
$ws.Cells.Item(3, 14).Value = 0.312473
$ws.Cells.Item(3, 15).Value = 0.0406458950147437
$ws.Cells.Item(3, 16).Value = 0.04064589501474371
$ws.Cells.Item(3, 17).Value = 0.01355407188255555
$ws.Cells.Item(3, 18).Value = 0.121986646943
$ws.Cells.Item(3, 19).Value = 0.01213737999623728
$ws.Cells.Item(3, 20).Value = 0.01213737999623728

# Row 4: ECs -> Inflammatory-Mac
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Bmp7"
$ws.Cells.Item(4, 3).Value = "Acvr2b"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.1301303333333333
$ws.Cells.Item(4, 8).Value = 0.390391
$ws.Cells.Item(4, 9).Value = 0.2986126887311651
$ws.Cells.Item(4, 10).Value = 0.2986126887311651
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.3139526666666667
$ws.Cells.Item(4, 14).Value = 0.9418580000000001
$ws.Cells.Item(4, 15).Value = 0.1225151017425393
$ws.Cells.Item(4, 16).Value = 0.1225151017425393
$ws.Cells.Item(4, 17).Value = 0.04085476516422223
$ws.Cells.Item(4, 18).Value = 0.367692886478
$ws.Cells.Item(4, 19).Value = 0.03658456394151191
$ws.Cells.Item(4, 20).Value = 0.03658456394151191

# Row 5: ECs -> MuSCs
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Bmp7"
$ws.Cells.Item(5, 3).Value = "Acvr2b"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 0.3333333333333333
$ws.Cells.Item(5, 7).Value = 0.1301303333333333
$ws.Cells.Item(5, 8).Value = 0.390391
$ws.Cells.Item(5, 9).Value = 0.2986126887311651
$ws.Cells.Item(5, 10).Value = 0.2986126887311651
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.5829876666666666
$ws.Cells.Item(5, 14).Value = 1.748963
$ws.Cells.Item(5, 15).Value = 0.2275017888991087
$ws.Cells.Item(5, 16).Value = 0.2275017888991087
$ws.Cells.Item(5, 17).Value = 0.07586437939255554
$ws.Cells.Item(5, 18).Value = 0.682779414533
$ws.Cells.Item(5, 19).Value = 0.06793492087431276
$ws.Cells.Item(5, 20).Value = 0.06793492087431278

# Row 6: ECs -> Neutrophils
$ws.Cells.Item(6, 1).Value = "ECs"
$ws.Cells.Item(6, 2).Value = "Bmp7"
$ws.Cells.Item(6, 3).Value = "Acvr2b"
$ws.Cells.Item(6, 4).Value = "Neutrophils"
$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 0.3333333333333333
$ws.Cells.Item(6, 7).Value = 0.1301303333333333
$ws.Cells.Item(6, 8).Value = 0.390391
$ws.Cells.Item(6, 9).Value = 0.2986126887311651
$ws.Cells.Item(6, 10).Value = 0.2986126887311651
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.3019996666666667
$ws.Cells.Item(6, 14).Value = 0.905999
$ws.Cells.Item(6, 15).Value = 0.117850631054404
$ws.Cells.Item(6, 16).Value = 0.117850631054404
$ws.Cells.Item(6, 17).Value = 0.03929931728988888
$ws.Cells.Item(6, 18).Value = 0.353693855609
$ws.Cells.Item(6, 19).Value = 0.03519169380782012
$ws.Cells.Item(6, 20).Value = 0.03519169380782012

# Row 7: ECs -> Resolving-Mac
$ws.Cells.Item(7, 1).Value = "ECs"
$ws.Cells.Item(7, 2).Value = "Bmp7"
$ws.Cells.Item(7, 3).Value = "Acvr2b"
$ws.Cells.Item(7, 4).Value = "Resolving-Mac"
$ws.Cells.Item(7, 5).Value = 1
$ws.Cells.Item(7, 6).Value = 0.3333333333333333
$ws.Cells.Item(7, 7).Value = 0.1301303333333333
$ws.Cells.Item(7, 8).Value = 0.390391
$ws.Cells.Item(7, 9).Value = 0.2986126887311651
$ws.Cells.Item(7, 10).Value = 0.2986126887311651
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.2807033333333333
$ws.Cells.Item(7, 14).Value = 0.8421099999999999
$ws.Cells.Item(7, 15).Value = 0.1095400711449175
$ws.Cells.Item(7, 16).Value = 0.1095400711449176
$ws.Cells.Item(7, 17).Value = 0.03652801833444444
$ws.Cells.Item(7, 18).Value = 0.3287521650099999
$ws.Cells.Item(7, 19).Value = 0.03271005516838694
$ws.Cells.Item(7, 20).Value = 0.03271005516838695

# Row 8: FAPs -> ECs
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Bmp7"
$ws.Cells.Item(8, 3).Value = "Acvr2b"
$ws.Cells.Item(8, 4).Value = "ECs"
$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 0.3333333333333333
$ws.Cells.Item(8, 7).Value = 0.2842186666666667
$ws.Cells.Item(8, 8).Value = 0.852656
$ws.Cells.Item(8, 9).Value = 0.6522022811047395
$ws.Cells.Item(8, 10).Value = 0.6522022811047393
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 0.978762
$ws.Cells.Item(8, 14).Value = 2.936286
$ws.Cells.Item(8, 15).Value = 0.3819465121442868
$ws.Cells.Item(8, 16).Value = 0.3819465121442868
$ws.Cells.Item(8, 17).Value = 0.278182430624
$ws.Cells.Item(8, 18).Value = 2.503641875616
$ws.Cells.Item(8, 19).Value = 0.249106386480503
$ws.Cells.Item(8, 20).Value = 0.2491063864805029

# Row 9: FAPs -> FAPs
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Bmp7"
$ws.Cells.Item(9, 3).Value = "Acvr2b"
$ws.Cells.Item(9, 4).Value = "FAPs"
$ws.Cells.Item(9, 5).Value = 1
$ws.Cells.Item(9, 6).Value = 0.3333333333333333
$ws.Cells.Item(9, 7).Value = 0.2842186666666667
$ws.Cells.Item(9, 8).Value = 0.852656
$ws.Cells.Item(9, 9).Value = 0.6522022811047395
$ws.Cells.Item(9, 10).Value = 0.6522022811047393
$ws.Cells.Item(9, 11).Value = 2
$ws.Cells.Item(9, 12).Value = 0.6666666666666666
$ws.Cells.Item(9, 13).Value = 0.1041576666666667
$ws.Cells.Item(9, 14).Value = 0.312473
$ws.Cells.Item(9, 15).Value = 0.0406458950147437
$ws.Cells.Item(9, 16).Value = 0.04064589501474371
$ws.Cells.Item(9, 17).Value = 0.02960355314311111
$ws.Cells.Item(9, 18).Value = 0.266431978288
$ws.Cells.Item(9, 19).Value = 0.0265093454461596
$ws.Cells.Item(9, 20).Value = 0.0265093454461596

# Row 10: FAPs -> Inflammatory-Mac
$ws.Cells.Item(10, 1).Value = "FAPs"
$ws.Cells.Item(10, 2).Value = "Bmp7"
$ws.Cells.Item(10, 3).Value = "Acvr2b"
$ws.Cells.Item(10, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(10, 5).Value = 1
$ws.Cells.Item(10, 6).Value = 0.3333333333333333
$ws.Cells.Item(10, 7).Value = 0.2842186666666667
$ws.Cells.Item(10, 8).Value = 0.852656
$ws.Cells.Item(10, 9).Value = 0.6522022811047395
$ws.Cells.Item(10, 10).Value = 0.6522022811047393
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.3139526666666667
$ws.Cells.Item(10, 14).Value = 0.9418580000000001
$ws.Cells.Item(10, 15).Value = 0.1225151017425393
$ws.Cells.Item(10, 16).Value = 0.1225151017425393
$ws.Cells.Item(10, 17).Value = 0.08923120831644446
$ws.Cells.Item(10, 18).Value = 0.8030808748480001
$ws.Cells.Item(10, 19).Value = 0.07990462882626338
$ws.Cells.Item(10, 20).Value = 0.07990462882626335

# Row 11: FAPs -> MuSCs
$ws.Cells.Item(11, 1).Value = "FAPs"
$ws.Cells.Item(11, 2).Value = "Bmp7"
$ws.Cells.Item(11, 3).Value = "Acvr2b"
$ws.Cells.Item(11, 4).Value = "MuSCs"
$ws.Cells.Item(11, 5).Value = 1
$ws.Cells.Item(11, 6).Value = 0.3333333333333333
$ws.Cells.Item(11, 7).Value = 0.2842186666666667
$ws.Cells.Item(11, 8).Value = 0.852656
$ws.Cells.Item(11, 9).Value = 0.6522022811047395
$ws.Cells.Item(11, 10).Value = 0.6522022811047393
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.5829876666666666
$ws.Cells.Item(11, 14).Value = 1.748963
$ws.Cells.Item(11, 15).Value = 0.2275017888991087
$ws.Cells.Item(11, 16).Value = 0.2275017888991087
$ws.Cells.Item(11, 17).Value = 0.1656959773031111
$ws.Cells.Item(11, 18).Value = 1.491263795728
$ws.Cells.Item(11, 19).Value = 0.1483771856754076
$ws.Cells.Item(11, 20).Value = 0.1483771856754076

# Row 12: FAPs -> Neutrophils
$ws.Cells.Item(12, 1).Value = "FAPs"
$ws.Cells.Item(12, 2).Value = "Bmp7"
$ws.Cells.Item(12, 3).Value = "Acvr2b"
$ws.Cells.Item(12, 4).Value = "Neutrophils"
$ws.Cells.Item(12, 5).Value = 1
$ws.Cells.Item(12, 6).Value = 0.3333333333333333
$ws.Cells.Item(12, 7).Value = 0.2842186666666667
$ws.Cells.Item(12, 8).Value = 0.852656
$ws.Cells.Item(12, 9).Value = 0.6522022811047395
$ws.Cells.Item(12, 10).Value = 0.6522022811047393
$ws.Cells.Item(12, 11).Value = 2
$ws.Cells.Item(12, 12).Value = 0.6666666666666666
$ws.Cells.Item(12, 13).Value = 0.3019996666666667
$ws.Cells.Item(12, 14).Value = 0.905999
$ws.Cells.Item(12, 15).Value = 0.117850631054404
$ws.Cells.Item(12, 16).Value = 0.117850631054404
$ws.Cells.Item(12, 17).Value = 0.08583394259377777
$ws.Cells.Item(12, 18).Value = 0.7725054833439999
$ws.Cells.Item(12, 19).Value = 0.07686245040331534
$ws.Cells.Item(12, 20).Value = 0.07686245040331532

# Row 13: FAPs -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = "FAPs"
$ws.Cells.Item(13, 2).Value = "Bmp7"
$ws.Cells.Item(13, 3).Value = "Acvr2b"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 1
$ws.Cells.Item(13, 6).Value = 0.3333333333333333
$ws.Cells.Item(13, 7).Value = 0.2842186666666667
$ws.Cells.Item(13, 8).Value = 0.852656
$ws.Cells.Item(13, 9).Value = 0.6522022811047395
$ws.Cells.Item(13, 10).Value = 0.6522022811047393
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.2807033333333333
$ws.Cells.Item(13, 14).Value = 0.8421099999999999
$ws.Cells.Item(13, 15).Value = 0.1095400711449175
$ws.Cells.Item(13, 16).Value = 0.1095400711449176
$ws.Cells.Item(13, 17).Value = 0.07978112712888888
$ws.Cells.Item(13, 18).Value = 0.7180301441599999
$ws.Cells.Item(13, 19).Value = 0.07144228427309068
$ws.Cells.Item(13, 20).Value = 0.07144228427309066

# Row 14: MuSCs -> ECs
$ws.Cells.Item(14, 1).Value = "MuSCs"
$ws.Cells.Item(14, 2).Value = "Bmp7"
$ws.Cells.Item(14, 3).Value = "Acvr2b"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.021434
$ws.Cells.Item(14, 8).Value = 0.064302
$ws.Cells.Item(14, 9).Value = 0.04918503016409543
$ws.Cells.Item(14, 10).Value = 0.04918503016409543
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 0.978762
$ws.Cells.Item(14, 14).Value = 2.936286
$ws.Cells.Item(14, 15).Value = 0.3819465121442868
$ws.Cells.Item(14, 16).Value = 0.3819465121442868
$ws.Cells.Item(14, 17).Value = 0.020978784708
$ws.Cells.Item(14, 18).Value = 0.188809062372
$ws.Cells.Item(14, 19).Value = 0.01878605072088779
$ws.Cells.Item(14, 20).Value = 0.01878605072088779

# Row 15: MuSCs -> FAPs
$ws.Cells.Item(15, 1).Value = "MuSCs"
$ws.Cells.Item(15, 2).Value = "Bmp7"
$ws.Cells.Item(15, 3).Value = "Acvr2b"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.021434
$ws.Cells.Item(15, 8).Value = 0.064302
$ws.Cells.Item(15, 9).Value = 0.04918503016409543
$ws.Cells.Item(15, 10).Value = 0.04918503016409543
$ws.Cells.Item(15, 11).Value = 2
$ws.Cells.Item(15, 12).Value = 0.6666666666666666
$ws.Cells.Item(15, 13).Value = 0.1041576666666667
$ws.Cells.Item(15, 14).Value = 0.312473
$ws.Cells.Item(15, 15).Value = 0.0406458950147437
$ws.Cells.Item(15, 16).Value = 0.04064589501474371
$ws.Cells.Item(15, 17).Value = 0.002232515427333333
$ws.Cells.Item(15, 18).Value = 0.020092638846
$ws.Cells.Item(15, 19).Value = 0.001999169572346825
$ws.Cells.Item(15, 20).Value = 0.001999169572346825

# Row 16: MuSCs -> Inflammatory-Mac
$ws.Cells.Item(16, 1).Value = "MuSCs"
$ws.Cells.Item(16, 2).Value = "Bmp7"
$ws.Cells.Item(16, 3).Value = "Acvr2b"
$ws.Cells.Item(16, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.021434
$ws.Cells.Item(16, 8).Value = 0.064302
$ws.Cells.Item(16, 9).Value = 0.04918503016409543
$ws.Cells.Item(16, 10).Value = 0.04918503016409543
$ws.Cells.Item(16, 11).Value = 2
$ws.Cells.Item(16, 12).Value = 0.6666666666666666
$ws.Cells.Item(16, 13).Value = 0.3139526666666667
$ws.Cells.Item(16, 14).Value = 0.9418580000000001
$ws.Cells.Item(16, 15).Value = 0.1225151017425393
$ws.Cells.Item(16, 16).Value = 0.1225151017425393
$ws.Cells.Item(16, 17).Value = 0.006729261457333334
$ws.Cells.Item(16, 18).Value = 0.060563353116
$ws.Cells.Item(16, 19).Value = 0.006025908974764016
$ws.Cells.Item(16, 20).Value = 0.006025908974764016

# Row 17: MuSCs -> MuSCs
$ws.Cells.Item(17, 1).Value = "MuSCs"
$ws.Cells.Item(17, 2).Value = "Bmp7"
$ws.Cells.Item(17, 3).Value = "Acvr2b"
$ws.Cells.Item(17, 4).Value = "MuSCs"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.021434
$ws.Cells.Item(17, 8).Value = 0.064302
$ws.Cells.Item(17, 9).Value = 0.04918503016409543
$ws.Cells.Item(17, 10).Value = 0.04918503016409543
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.5829876666666666
$ws.Cells.Item(17, 14).Value = 1.748963
$ws.Cells.Item(17, 15).Value = 0.2275017888991087
$ws.Cells.Item(17, 16).Value = 0.2275017888991087
$ws.Cells.Item(17, 17).Value = 0.01249575764733333
$ws.Cells.Item(17, 18).Value = 0.112461818826
$ws.Cells.Item(17, 19).Value = 0.01118968234938833
$ws.Cells.Item(17, 20).Value = 0.01118968234938833

# Row 18: MuSCs -> Neutrophils
$ws.Cells.Item(18, 1).Value = "MuSCs"
$ws.Cells.Item(18, 2).Value = "Bmp7"
$ws.Cells.Item(18, 3).Value = "Acvr2b"
$ws.Cells.Item(18, 4).Value = "Neutrophils"
$ws.Cells.Item(18, 5).Value = 1
$ws.Cells.Item(18, 6).Value = 0.3333333333333333
$ws.Cells.Item(18, 7).Value = 0.021434
$ws.Cells.Item(18, 8).Value = 0.064302
$ws.Cells.Item(18, 9).Value = 0.04918503016409543
$ws.Cells.Item(18, 10).Value = 0.04918503016409543
$ws.Cells.Item(18, 11).Value = 2
$ws.Cells.Item(18, 12).Value = 0.6666666666666666
$ws.Cells.Item(18, 13).Value = 0.3019996666666667
$ws.Cells.Item(18, 14).Value = 0.905999
$ws.Cells.Item(18, 15).Value = 0.117850631054404
$ws.Cells.Item(18, 16).Value = 0.117850631054404
$ws.Cells.Item(18, 17).Value = 0.006473060855333333
$ws.Cells.Item(18, 18).Value = 0.058257547698
$ws.Cells.Item(18, 19).Value = 0.005796486843268542
$ws.Cells.Item(18, 20).Value = 0.005796486843268543

# Row 19: MuSCs -> Resolving-Mac
$ws.Cells.Item(19, 1).Value = "MuSCs"
$ws.Cells.Item(19, 2).Value = "Bmp7"
$ws.Cells.Item(19, 3).Value = "Acvr2b"
$ws.Cells.Item(19, 4).Value = "Resolving-Mac"
$ws.Cells.Item(19, 5).Value = 1
$ws.Cells.Item(19, 6).Value = 0.3333333333333333
$ws.Cells.Item(19, 7).Value = 0.021434
$ws.Cells.Item(19, 8).Value = 0.064302
$ws.Cells.Item(19, 9).Value = 0.04918503016409543
$ws.Cells.Item(19, 10).Value = 0.04918503016409543
$ws.Cells.Item(19, 11).Value = 3
$ws.Cells.Item(19, 12).Value = 1
$ws.Cells.Item(19, 13).Value = 0.2807033333333333
$ws.Cells.Item(19, 14).Value = 0.8421099999999999
$ws.Cells.Item(19, 15).Value = 0.1095400711449175
$ws.Cells.Item(19, 16).Value = 0.1095400711449176
$ws.Cells.Item(19, 17).Value = 0.006016595246666666
$ws.Cells.Item(19, 18).Value = 0.05414935721999999
$ws.Cells.Item(19, 19).Value = 0.005387731703439928
$ws.Cells.Item(19, 20).Value = 0.005387731703439929

